$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 1.68
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("W2").Value = 11
$ws.Range("AC2").Value = 7
$ws.Range("AQ2").Value = 126
$ws.Range("AS2").Value = 401
$ws.Range("AU2").Value = 9.5
$ws.Range("BA2").Value = 67
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("BD3").Value = 126
$ws.Range("BD4").Value = 151
$ws.Range("BD5").Value = 126
$ws.Range("G6").Value = 2.25
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 3.3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 2.05
$ws.Range("L6").Value = 4
$ws.Range("X6").Value = 10
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 21
$ws.Range("AC6").Value = 8
$ws.Range("AH6").Value = 9
$ws.Range("AL6").Value = 29
$ws.Range("AO6").Value = 13
$ws.Range("AQ6").Value = 41
$ws.Range("AR6").Value = 67
$ws.Range("AV6").Value = 51
$ws.Range("AW6").Value = 5
$ws.Range("AX6").Value = 19
$ws.Range("AZ6").Value = 67
$ws.Range("G7").Value = 2.05
$ws.Range("L8").Value = 3.1
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 1.48
$ws.Range("AF8").Value = 81
$ws.Range("AI8").Value = 9.5
$ws.Range("AN8").Value = 5
$ws.Range("AV8").Value = 81
$ws.Range("AZ8").Value = 51
$ws.Range("G9").Value = 2.45
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 3.25
$ws.Range("X9").Value = 11
$ws.Range("AH9").Value = 7.5
$ws.Range("AI9").Value = 13
$ws.Range("AW9").Value = 4.75
$ws.Range("AY9").Value = 29
$ws.Range("AZ9").Value = 51
$ws.Range("BB9").Value = 251
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 8.5
$ws.Range("Q10").Value = 1.97
$ws.Range("R10").Value = 1.77
$ws.Range("AA10").Value = 13
$ws.Range("AJ10").Value = 34
$ws.Range("AK10").Value = 151
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("S11").Value = 1.5
$ws.Range("T11").Value = 2.5
$ws.Range("AA11").Value = 17
$ws.Range("AD11").Value = 7
$ws.Range("AH11").Value = 11
$ws.Range("AN11").Value = 3.5
$ws.Range("AT11").Value = 2.5
$ws.Range("H12").Value = 4.33
$ws.Range("I12").Value = 8
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.9
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.75
$ws.Range("W12").Value = 6
$ws.Range("Y12").Value = 8.5
$ws.Range("Z12").Value = 9.5
$ws.Range("AB12").Value = 29
$ws.Range("AG12").Value = 1000
$ws.Range("AR12").Value = 51
$ws.Range("AT12").Value = 2.75
$ws.Range("BC13").Value = 126
$ws.Range("BD13").Value = 126
$ws.Range("G14").Value = 1.75
$ws.Range("I14").Value = 4.5
$ws.Range("J14").Value = 2.38
$ws.Range("L14").Value = 4.75
$ws.Range("U14").Value = 1.8
$ws.Range("V14").Value = 1.91
$ws.Range("X14").Value = 8.5
$ws.Range("AA14").Value = 13
$ws.Range("AC14").Value = 11
$ws.Range("AI14").Value = 23
$ws.Range("AJ14").Value = 15
$ws.Range("AK14").Value = 51
$ws.Range("AO14").Value = 9
$ws.Range("AW14").Value = 6.5
$ws.Range("BA14").Value = 101
$ws.Range("BC14").Value = 151
$ws.Range("M16").Value = 1.04
$ws.Range("N16").Value = 13
$ws.Range("O16").Value = 1.22
$ws.Range("P16").Value = 4.33
$ws.Range("Q16").Value = 1.73
$ws.Range("R16").Value = 2.1
$ws.Range("O17").Value = 1.4
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 2.2
$ws.Range("R17").Value = 1.67
$ws.Range("G18").Value = 1.4
$ws.Range("H18").Value = 4.5
$ws.Range("I18").Value = 8.5
$ws.Range("J18").Value = 1.91
$ws.Range("K18").Value = 2.4
$ws.Range("O18").Value = 1.22
$ws.Range("P18").Value = 4
$ws.Range("Q18").Value = 1.8
$ws.Range("R18").Value = 2
$ws.Range("S18").Value = 1.33
$ws.Range("T18").Value = 3.25
$ws.Range("W18").Value = 6.5
$ws.Range("Y18").Value = 8.5
$ws.Range("AH18").Value = 19
$ws.Range("AJ18").Value = 23
$ws.Range("AQ18").Value = 19
$ws.Range("AT18").Value = 3.25
$ws.Range("AV18").Value = 51
$ws.Range("BB18").Value = 301
$ws.Range("G21").Value = 1.85
$ws.Range("H21").Value = 3.7
$ws.Range("I21").Value = 4
$ws.Range("J21").Value = 2.4
$ws.Range("K21").Value = 2.4
$ws.Range("O21").Value = 1.18
$ws.Range("P21").Value = 4.5
$ws.Range("Q21").Value = 1.62
$ws.Range("R21").Value = 2.25
$ws.Range("S21").Value = 1.29
$ws.Range("T21").Value = 3.5
$ws.Range("W21").Value = 10
$ws.Range("AI21").Value = 23
$ws.Range("AL21").Value = 29
$ws.Range("AO21").Value = 9.5
$ws.Range("AT21").Value = 3.5
$ws.Range("BC21").Value = 351
